{"js": "// Apply the four resume text updates (German -> English-with-company-name\n// style date lines, plus a \"Telefon:\" prefix on the phone number), as\n// described by the diff. Each target string is unique in the document and\n// lives in its own single-run paragraph, so a straightforward\n// search + insertText(Replace) per pair is sufficient and keeps the\n// surrounding run formatting (rPr) untouched.\n\nconst replacements = [\n  {\n    find: \"(123) 456-7890\",\n    replace: \"Telefon: (123) 456-7890\",\n  },\n  {\n    find: \"Lead Animator (Januar 2018\\u2013Heute)\",\n    replace: \"ABC Studios: Lead Animator (Jan 2018 - Present)\",\n  },\n  {\n    find: \"Senior Animator (Juni 2015\\u2013Dezember 2017)\",\n    replace: \"XYZ Media: Senior Animator (Jun 2015 - Dez 2017)\",\n  },\n  {\n    find: \"Junior Animator (September 2012\\u2013Mai 2015)\",\n    replace: \"MNO Entertainment: Junior Animator (Sep 2012 - Mai 2015)\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"(123) 456-7890\"; Replace = \"Telefon: (123) 456-7890\" },\n    @{ Find = \"Lead Animator (Januar 2018\u2013Heute)\"; Replace = \"ABC Studios: Lead Animator (Jan 2018 - Present)\" },\n    @{ Find = \"Senior Animator (Juni 2015\u2013Dezember 2017)\"; Replace = \"XYZ Media: Senior Animator (Jun 2015 - Dez 2017)\" },\n    @{ Find = \"Junior Animator (September 2012\u2013Mai 2015)\"; Replace = \"MNO Entertainment: Junior Animator (Sep 2012 - Mai 2015)\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($r.Find, $true, $true, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $($r.Find)\"\n    }\n}\n"}
